# g20.3 -> g20.4 : atualização da fonte
# Refresh the data source: new date (31/12/2010), updated IVS values,
# and re-ordered / renamed regions. Also restyle the header row with a
# thin border around each header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, in final row order (row -> Região, Ano, Valor, Colocação)
$data = @(
    @{ Row = 2;  Regiao = "Maranhão";   Ano = "31/12/2010"; Valor = 0.521; Colocacao = "1º" },
    @{ Row = 3;  Regiao = "Amazonas";   Ano = "31/12/2010"; Valor = 0.488; Colocacao = "2º" },
    @{ Row = 4;  Regiao = "Pará";       Ano = "31/12/2010"; Valor = 0.469; Colocacao = "3º" },
    @{ Row = 5;  Regiao = "Alagoas";    Ano = "31/12/2010"; Valor = 0.461; Colocacao = "4º" },
    @{ Row = 6;  Regiao = "Acre";       Ano = "31/12/2010"; Valor = 0.443; Colocacao = "5º" },
    @{ Row = 7;  Regiao = "Pernambuco"; Ano = "31/12/2010"; Valor = 0.414; Colocacao = "6º" },
    @{ Row = 8;  Regiao = "Sergipe";    Ano = "31/12/2010"; Valor = 0.393; Colocacao = "10º" },
    @{ Row = 9;  Regiao = "Brasil";     Ano = "31/12/2010"; Valor = 0.326; Colocacao = "" },
    @{ Row = 10; Regiao = "Nordeste";   Ano = "31/12/2010"; Valor = 0.408; Colocacao = "" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Regiao
    $ws.Cells.Item($r, 3).Value = $item.Ano
    $ws.Cells.Item($r, 4).Value = $item.Valor
    $ws.Cells.Item($r, 5).Value = $item.Colocacao
}

# Add a thin border all around the header row cells (A1:E1) and align
# the header text centered/top, matching the refreshed style sheet.
$header = $ws.Range("A1:E1")
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Reset page margins to Excel's own defaults (points: 1 inch = 72pt).
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

Write-Output "done"
